$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above row 459, shifting existing rows 459-478 down to 465-484
$ws.Range("A459:T464").EntireRow.Insert()

# Row 459
$ws.Range("A459").Value = 7
$ws.Range("B459").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C459").Value = 'Ñuble'
$ws.Range("D459").Value = 45147
$ws.Range("E459").Value = 16
$ws.Range("F459").Value = 'Fruta'
$ws.Range("G459").Value = 100104
$ws.Range("H459").Value = 'Frutos de pepita'
$ws.Range("I459").Value = 100104005
$ws.Range("J459").Value = 'Pera'
$ws.Range("K459").Value = 'Forelle'
$ws.Range("L459").Value = 'Especial'
$ws.Range("M459").Value = 60
$ws.Range("N459").Value = 12000
$ws.Range("O459").Value = 12000
$ws.Range("P459").Value = 12000
$ws.Range("Q459").Value = '$/bandeja 18 kilos granel'
$ws.Range("R459").Value = 'Región de O''Higgins'
$ws.Range("S459").Value = 667
$ws.Range("T459").Value = 18

# Row 460
$ws.Range("A460").Value = 7
$ws.Range("B460").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C460").Value = 'Ñuble'
$ws.Range("D460").Value = 45147
$ws.Range("E460").Value = 16
$ws.Range("F460").Value = 'Fruta'
$ws.Range("G460").Value = 100104
$ws.Range("H460").Value = 'Frutos de pepita'
$ws.Range("I460").Value = 100104005
$ws.Range("J460").Value = 'Pera'
$ws.Range("K460").Value = 'Forelle'
$ws.Range("L460").Value = 'Primera'
$ws.Range("M460").Value = 60
$ws.Range("N460").Value = 10000
$ws.Range("O460").Value = 10000
$ws.Range("P460").Value = 10000
$ws.Range("Q460").Value = '$/bandeja 18 kilos granel'
$ws.Range("R460").Value = 'Región de O''Higgins'
$ws.Range("S460").Value = 556
$ws.Range("T460").Value = 18

# Row 461
$ws.Range("A461").Value = 7
$ws.Range("B461").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C461").Value = 'Ñuble'
$ws.Range("D461").Value = 45147
$ws.Range("E461").Value = 16
$ws.Range("F461").Value = 'Fruta'
$ws.Range("G461").Value = 100104
$ws.Range("H461").Value = 'Frutos de pepita'
$ws.Range("I461").Value = 100104005
$ws.Range("J461").Value = 'Pera'
$ws.Range("K461").Value = 'Packham''s Triumph'
$ws.Range("L461").Value = 'Especial'
$ws.Range("M461").Value = 60
$ws.Range("N461").Value = 12000
$ws.Range("O461").Value = 12000
$ws.Range("P461").Value = 12000
$ws.Range("Q461").Value = '$/bandeja 18 kilos granel'
$ws.Range("R461").Value = 'Región de O''Higgins'
$ws.Range("S461").Value = 667
$ws.Range("T461").Value = 18

# Row 462
$ws.Range("A462").Value = 7
$ws.Range("B462").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C462").Value = 'Ñuble'
$ws.Range("D462").Value = 45147
$ws.Range("E462").Value = 16
$ws.Range("F462").Value = 'Fruta'
$ws.Range("G462").Value = 100104
$ws.Range("H462").Value = 'Frutos de pepita'
$ws.Range("I462").Value = 100104005
$ws.Range("J462").Value = 'Pera'
$ws.Range("K462").Value = 'Packham''s Triumph'
$ws.Range("L462").Value = 'Primera'
$ws.Range("M462").Value = 60
$ws.Range("N462").Value = 10000
$ws.Range("O462").Value = 10000
$ws.Range("P462").Value = 10000
$ws.Range("Q462").Value = '$/bandeja 18 kilos granel'
$ws.Range("R462").Value = 'Región de O''Higgins'
$ws.Range("S462").Value = 556
$ws.Range("T462").Value = 18

# Row 463
$ws.Range("A463").Value = 7
$ws.Range("B463").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C463").Value = 'Ñuble'
$ws.Range("D463").Value = 45147
$ws.Range("E463").Value = 16
$ws.Range("F463").Value = 'Fruta'
$ws.Range("G463").Value = 100104
$ws.Range("H463").Value = 'Frutos de pepita'
$ws.Range("I463").Value = 100104005
$ws.Range("J463").Value = 'Pera'
$ws.Range("K463").Value = 'Winter Nelis'
$ws.Range("L463").Value = 'Primera'
$ws.Range("M463").Value = 60
$ws.Range("N463").Value = 10000
$ws.Range("O463").Value = 10000
$ws.Range("P463").Value = 10000
$ws.Range("Q463").Value = '$/bandeja 18 kilos granel'
$ws.Range("R463").Value = 'Región de O''Higgins'
$ws.Range("S463").Value = 556
$ws.Range("T463").Value = 18

# Row 464
$ws.Range("A464").Value = 7
$ws.Range("B464").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C464").Value = 'Ñuble'
$ws.Range("D464").Value = 45147
$ws.Range("E464").Value = 16
$ws.Range("F464").Value = 'Fruta'
$ws.Range("G464").Value = 100104
$ws.Range("H464").Value = 'Frutos de pepita'
$ws.Range("I464").Value = 100104005
$ws.Range("J464").Value = 'Pera'
$ws.Range("K464").Value = 'Winter Nelis'
$ws.Range("L464").Value = 'Segunda'
$ws.Range("M464").Value = 60
$ws.Range("N464").Value = 8000
$ws.Range("O464").Value = 8000
$ws.Range("P464").Value = 8000
$ws.Range("Q464").Value = '$/bandeja 18 kilos granel'
$ws.Range("R464").Value = 'Región de O''Higgins'
$ws.Range("S464").Value = 444
$ws.Range("T464").Value = 18
